# Add a "best guess" column (D) with a verification-dataset formula,
# mirroring the existing pattern already used for column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("D1").Value = "best guess"

# Row 2 is a standalone (non-shared) formula, same convention as C2.
$ws.Range("D2").Formula = "=(0.3302*SIN(2.043*B2)*TANH(1.2639*A2))"

# Rows 3-66 share one formula (mirrors C3:C66's shared group).
$ws.Range("D3:D66").Formula = "=(0.3302*SIN(2.043*B3)*TANH(1.2639*A3))"

# Rows 67-82 share another formula (mirrors C67:C74's shared group).
$ws.Range("D67:D82").Formula = "=(0.3302*SIN(2.043*B67)*TANH(1.2639*A67))"

# Column D width, as close as the engine's rounding allows to 10.7109375.
$ws.Range("D1").ColumnWidth = 9.8

# Selection moves to S9 in the saved file.
$ws.Range("S9").Select() | Out-Null
